$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("INDEX (DO NOT MODIFY)")
$ws.Columns.Item(1).Insert()

# Copy formatting from column B (the original column A) into the new column A
$ws.Range("B1:B7").Copy()
$ws.Range("A1:A7").PasteSpecial(-4122)
$ws.Columns.Item(1).ColumnWidth = 22.14

# Header + index values for the new column A
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 126
$ws.Range("A3").Value = 127
$ws.Range("A4").Value = 128
$ws.Range("A5").Value = 169
$ws.Range("A6").Value = 206
$ws.Range("A7").Value = 289

# Uppercase / rewording of several header labels (now shifted one column right)
$ws.Range("C1").Value = "REGION"
$ws.Range("R1").Value = " TARGET COMPLETION DATE "
$ws.Range("S1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("T1").Value = "PROJECT ID"
$ws.Range("U1").Value = "CONTRACT ID"
$ws.Range("V1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("W1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("X1").Value = "BID OPENING"
$ws.Range("Y1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("Z1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("AA1").Value = "NAME OF CONTRACTOR"
$ws.Range("AB1").Value = "OTHER REMARKS"

Write-Host "done"
